$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New helper column G: a constant "1" marker on every existing data row
# (the commit adds a hidden sort-order helper column so the sheet can be
# re-sorted in either direction later).
for ($r = 1; $r -le 94; $r++) {
    $ws.Cells.Item($r, 7).Value = 1
}

# --- Four new module rows (95-98): German translations added at the end
# of the list, each also gets the new column G marker.
$newRows = @(
    @("German Button",          "BigButtonTranslated",    1, "./resources/modules/German Button.pdf",          "Malde, Tharagon", "2017-01-25"),
    @("German Morse Code",      "MorseCodeTranslated",    2, "./resources/modules/German Morse Code.pdf",      "Malde, Tharagon", "2017-01-25"),
    @("German Password",        "PasswordsTranslated",    2, "./resources/modules/German Password.pdf",        "Malde, Tharagon", "2017-01-25"),
    @("German Who's On First",  "WhosOnFirstTranslated",  2, "./resources/modules/German Who's On First.pdf",  "Malde, Tharagon", "2017-01-25")
)

$r = 95
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = 1
    $r = $r + 1
}

# --- Column width updates: column B widened (longer translated-module
# names), new column G sized like the other thin spacer columns (C).
$ws.Range("B1").ColumnWidth = 31.833333333333332
$ws.Range("G1").ColumnWidth = 1

# --- Selection / scroll position left by the editor after appending rows.
$excel.ActiveWindow.ScrollRow = 77
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K91").Select()
